$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a day-end "low stock" report. A new stock-out item
# (INFLUVAC TETRA ...) needs to be inserted as a new data row right after
# the "EGY TOTAVIT SYRUP 100 ML" row (which currently sits at row 9), before
# "MELLITOFIX 10MG 30 F.C. TABS" (currently row 10). Inserting the row
# pushes all the following data rows, the totals row and the footer row
# down by one. After that we fix up the running item counter, the grand
# total and the generated-at timestamp.
# ---------------------------------------------------------------------------

# Remove the merged ranges that live in/under the row we are about to
# insert into, so Excel doesn't choke on inserting inside a merge and so we
# can re-create them afterwards in the right row order.
$ws.Range("A10:Q18").UnMerge()

# Insert a new blank row at row 10 - shifts rows 10-18 down to 11-19.
$ws.Rows.Item(10).Insert()

# Give the new row the same per-cell styling as the row above it (row 9),
# restricted to the columns actually used by the table (A:Q) so we don't
# smear formatting across the whole 16384-column row.
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = 24.75

# Re-create the merged cells for the table rows (10-17), the totals row
# (18) and the footer row (19).
for ($r = 10; $r -le 17; $r++) {
    $ws.Range("A$r" + ":B$r").Merge()
    $ws.Range("C$r" + ":G$r").Merge()
    $ws.Range("H$r" + ":K$r").Merge()
    $ws.Range("L$r" + ":M$r").Merge()
    $ws.Range("N$r" + ":O$r").Merge()
}
$ws.Range("P18:Q18").Merge()
$ws.Range("A19:F19").Merge()
$ws.Range("G19:I19").Merge()
$ws.Range("K19:Q19").Merge()

# Fill in the new item's data in row 10.
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "INFLUVAC TETRA 1 PREF.SYRINGE SUSP. FOR INJ S.C. 0.5ML."
$ws.Range("H10").Value = "0:0"
$ws.Range("N10").Value = "300.00"
$ws.Range("P10").Value = "600.0000"
$ws.Range("Q10").Value = "2:0"

# Renumber the "م" (item index) column for the rows that shifted down.
for ($r = 11; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# Update the grand total (increases by the new row's sell-price, 600).
$ws.Range("P18").Value = 1047.3199999999999

# Update the generated-at timestamp in the footer.
$ws.Range("A19").Value = "Thursday, 18 September, 2025 11:36 AM"
